{"js": "// Replace each two-digit-division answer cell's text with the updated value.\n// Each \"before\" string below is unique in the document, so a targeted\n// search + insertText(\"Replace\") on every hit keeps run formatting\n// (font/size) untouched while only the text content changes.\nconst replacements = [\n  [\"36\u00f79=4, 0\", \"34\u00f79=3, 7\"],\n  [\"19\u00f76=3, 1\", \"53\u00f74=13, 1\"],\n  [\"51\u00f72=25, 1\", \"98\u00f75=19, 3\"],\n  [\"58\u00f79=6, 4\", \"89\u00f72=44, 1\"],\n  [\"67\u00f75=13, 2\", \"89\u00f73=29, 2\"],\n  [\"62\u00f77=8, 6\", \"18\u00f75=3, 3\"],\n  [\"46\u00f72=23, 0\", \"76\u00f76=12, 4\"],\n  [\"17\u00f74=4, 1\", \"30\u00f75=6, 0\"],\n  [\"12\u00f76=2, 0\", \"50\u00f79=5, 5\"],\n  [\"51\u00f79=5, 6\", \"71\u00f75=14, 1\"],\n  [\"25\u00f76=4, 1\", \"62\u00f75=12, 2\"],\n  [\"99\u00f74=24, 3\", \"83\u00f74=20, 3\"],\n  [\"68\u00f72=34, 0\", \"36\u00f77=5, 1\"],\n  [\"11\u00f76=1, 5\", \"70\u00f74=17, 2\"],\n  [\"48\u00f76=8, 0\", \"48\u00f72=24, 0\"],\n  [\"54\u00f78=6, 6\", \"77\u00f79=8, 5\"],\n  [\"59\u00f78=7, 3\", \"28\u00f74=7, 0\"],\n  [\"65\u00f79=7, 2\", \"59\u00f76=9, 5\"],\n  [\"30\u00f77=4, 2\", \"29\u00f78=3, 5\"],\n  [\"92\u00f76=15, 2\", \"89\u00f73=29, 2\"],\n  [\"17\u00f75=3, 2\", \"92\u00f72=46, 0\"],\n  [\"88\u00f79=9, 7\", \"73\u00f79=8, 1\"],\n  [\"70\u00f77=10, 0\", \"95\u00f74=23, 3\"],\n  [\"57\u00f75=11, 2\", \"99\u00f75=19, 4\"],\n  [\"87\u00f77=12, 3\", \"65\u00f73=21, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division answer cell's text with the updated value.\n# Each \"before\" string is unique in the document, so Find/Replace over the\n# whole document body for each pair retargets only the intended cell while\n# leaving run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('36\u00f79=4, 0', '34\u00f79=3, 7'),\n    @('19\u00f76=3, 1', '53\u00f74=13, 1'),\n    @('51\u00f72=25, 1', '98\u00f75=19, 3'),\n    @('58\u00f79=6, 4', '89\u00f72=44, 1'),\n    @('67\u00f75=13, 2', '89\u00f73=29, 2'),\n    @('62\u00f77=8, 6', '18\u00f75=3, 3'),\n    @('46\u00f72=23, 0', '76\u00f76=12, 4'),\n    @('17\u00f74=4, 1', '30\u00f75=6, 0'),\n    @('12\u00f76=2, 0', '50\u00f79=5, 5'),\n    @('51\u00f79=5, 6', '71\u00f75=14, 1'),\n    @('25\u00f76=4, 1', '62\u00f75=12, 2'),\n    @('99\u00f74=24, 3', '83\u00f74=20, 3'),\n    @('68\u00f72=34, 0', '36\u00f77=5, 1'),\n    @('11\u00f76=1, 5', '70\u00f74=17, 2'),\n    @('48\u00f76=8, 0', '48\u00f72=24, 0'),\n    @('54\u00f78=6, 6', '77\u00f79=8, 5'),\n    @('59\u00f78=7, 3', '28\u00f74=7, 0'),\n    @('65\u00f79=7, 2', '59\u00f76=9, 5'),\n    @('30\u00f77=4, 2', '29\u00f78=3, 5'),\n    @('92\u00f76=15, 2', '89\u00f73=29, 2'),\n    @('17\u00f75=3, 2', '92\u00f72=46, 0'),\n    @('88\u00f79=9, 7', '73\u00f79=8, 1'),\n    @('70\u00f77=10, 0', '95\u00f74=23, 3'),\n    @('57\u00f75=11, 2', '99\u00f75=19, 4'),\n    @('87\u00f77=12, 3', '65\u00f73=21, 2'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
